$d = $word.ActiveDocument

$pairs = @(
    @("42÷6=", "40÷2="),
    @("87÷8=", "34÷6="),
    @("13÷7=", "47÷3="),
    @("60÷9=", "53÷3="),
    @("17÷8=", "84÷4="),
    @("31÷2=", "32÷4="),
    @("29÷3=", "85÷3="),
    @("99÷8=", "65÷4="),
    @("54÷8=", "42÷3="),
    @("31÷6=", "46÷9="),
    @("92÷4=", "71÷9="),
    @("13÷9=", "45÷8="),
    @("83÷9=", "94÷6="),
    @("83÷3=", "11÷9="),
    @("93÷6=", "20÷5="),
    @("90÷2=", "69÷8="),
    @("12÷7=", "98÷2="),
    @("16÷7=", "68÷7="),
    @("24÷4=", "60÷4="),
    @("61÷5=", "17÷7="),
    @("51÷6=", "70÷2="),
    @("45÷3=", "49÷8="),
    @("89÷7=", "27÷2="),
    @("75÷5=", "60÷8="),
    @("49÷5=", "58÷6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
